$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-14
# from 2023-10-22 (45221) to 2023-10-25 (45224)
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
